$wb = $excel.ActiveWorkbook

# --- "Coupling Parameters" sheet edits -------------------------------------
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2089 -> 2060
$ws.Range("B4").Value = 2060

# Power plants year: turn the literal 2050 into a formula referencing the
# Start Year cell (value stays 2050 since B3 = 2050)
$ws.Range("B5").Formula = "=B3"

# capacity_remuneration_mechanism: was the combined description string in
# B45 with a long explanatory note in C45. Now just "capacity_market" is
# chosen via a dropdown, and the note cell is cleared.
$ws.Range("B45").Value = "capacity_market"
$ws.Range("C45").ClearContents()

# --- add a new "Sheet1" holding the allowed capacity-mechanism choices -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet1"
$newSheet.Range("A1").Value = "capacity_market"
$newSheet.Range("A2").Value = "capacity_subscription"
$newSheet.Range("A3").Value = "strategic_reserve_ger"
$newSheet.Range("A4").Value = "strategic_reserve_swe"
$newSheet.Range("A5").Value = "forward_capacity_market"

# --- data validation dropdown on B45 referencing the new list --------------
$ws.Range("B45").Validation.Add(3, 1, 1, '=Sheet1!$A$1:$A$5')
$ws.Range("B45").Validation.IgnoreBlank = $true
$ws.Range("B45").Validation.InCellDropdown = $true

# --- restore the active sheet / selection -----------------------------------
$ws.Activate()
$ws.Range("M5").Select() | Out-Null
